$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before the existing row 18),
# pushing the existing rows 18-100 down to become rows 20-102.
$ws.Rows.Item(18).EntireRow.Insert()
$ws.Rows.Item(18).EntireRow.Insert()

# Populate the first new row (row 18) with a new weekly record.
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value = 44802
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 100114002
$ws.Range("G18").Value = "Camote"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 520
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 14000
$ws.Range("N18").Value = "$/caja 18 kilos"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 778
$ws.Range("Q18").Value = 18
$ws.Range("R18").Value = "Hortaliza"

# Populate the second new row (row 19) with a new weekly record.
$ws.Range("A19").Value = 9
$ws.Range("B19").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44802
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100114002
$ws.Range("G19").Value = "Camote"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 790
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 14000
$ws.Range("N19").Value = "$/malla 18 kilos"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 778
$ws.Range("Q19").Value = 18
$ws.Range("R19").Value = "Hortaliza"

# Make sure the date column keeps the same date number format style as the
# rest of the column (style index carried from the row insert already gives
# this, but set explicitly to be safe).
$ws.Range("D18:D19").NumberFormat = $ws.Range("D20").NumberFormat
